# Started EMP removal + major/minor action revamp
# Add the new "Increased Accuracy" skill as a new row (63) on the skills sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New skill row: Name / SkillLevels / Effect / (no Prerequisite) / Include
$ws.Range("A63").Value = "Increased Accuracy"
$ws.Range("B63").Value = 1
$ws.Range("C63").Value = "When performing an accuracy check, you may add a bonus to the check equal to the amount that the spellcasting check exceeds the difficulty value. "
$ws.Range("E63").Value = 1

# Match the row height used by similarly-wrapped rows (e.g. row 62 above it).
$ws.Rows.Item(63).RowHeight = 23.85

# Move the selection to the newly added row, mirroring the authored view-state
# change (selection moved from C62 to D63 as work continued onto the new row).
$null = $ws.Range("D63").Select()
